$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.771.60"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "3.083.37"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.57"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.83"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.076.32"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +5.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.51"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "3.579.12"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "63.745.56"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "3.082.50"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.53"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.13"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.80"
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.22"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.36"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "56.90"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "496.13"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").Value = "3.275.46"
$ws.Range("E38").Value = "  +7.16%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0794"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.14"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").Value = "0.0₃0536"
$ws.Range("E47").Value = "  +6.87%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.80"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.84"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  -7.21%  "
